$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D/E data columns (rows 2-51) to Text format before writing, so that
# numeric-looking strings (e.g. "581.50", "0.998") are preserved verbatim as
# text rather than being auto-converted to numbers (which would drop
# trailing zeros / change representation).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '64.859.86'
$ws.Range("E2").Value = '  +3.48%  '
$ws.Range("D3").Value = '2.544.07'
$ws.Range("E3").Value = '  +3.33%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '581.50'
$ws.Range("D6").Value = '153.21'
$ws.Range("E6").Value = '  +4.10%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +1.13%  '
$ws.Range("D9").Value = '2.546.92'
$ws.Range("E9").Value = '  +3.45%  '
$ws.Range("E10").Value = '  +1.62%  '
$ws.Range("E11").Value = '  -1.65%  '
$ws.Range("E12").Value = '  +0.54%  '
$ws.Range("D13").Value = '0.356'
$ws.Range("E13").Value = '  +0.63%  '
$ws.Range("D14").Value = '29.34'
$ws.Range("E14").Value = '  +0.63%  '
$ws.Range("D15").Value = '0.0000180'
$ws.Range("E15").Value = '  +2.30%  '
$ws.Range("D16").Value = '3.006.78'
$ws.Range("E16").Value = '  +3.55%  '
$ws.Range("D17").Value = '64.830.05'
$ws.Range("E17").Value = '  +3.61%  '
$ws.Range("D18").Value = '2.551.96'
$ws.Range("E18").Value = '  +3.37%  '
$ws.Range("D19").Value = '8.07'
$ws.Range("E19").Value = '  +1.74%  '
$ws.Range("D20").Value = '11.03'
$ws.Range("E20").Value = '  +0.88%  '
$ws.Range("E21").Value = '  +3.67%  '
$ws.Range("D22").Value = '329.76'
$ws.Range("E22").Value = '  +1.35%  '
$ws.Range("D23").Value = '2.23'
$ws.Range("E23").Value = '  +2.30%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").Value = '10.12'
$ws.Range("E25").Value = '  +0.39%  '
$ws.Range("D26").Value = '65.86'
$ws.Range("E26").Value = '  +0.94%  '
$ws.Range("D27").Value = '637.44'
$ws.Range("E27").Value = '  -0.75%  '
$ws.Range("E28").Value = '  +7.41%  '
$ws.Range("D29").Value = '2.668.92'
$ws.Range("E29").Value = '  +3.37%  '
$ws.Range("E30").Value = '  +4.40%  '
$ws.Range("D31").Value = '0.998'
$ws.Range("E31").Value = '  -0.25%  '
$ws.Range("D32").Value = '8.09'
$ws.Range("E32").Value = '  +2.19%  '
$ws.Range("E33").Value = '  +2.69%  '
$ws.Range("D34").Value = '0.139'
$ws.Range("E34").Value = '  +4.06%  '
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("D36").Value = '1.58'
$ws.Range("E36").Value = '  +3.87%  '
$ws.Range("D37").Value = '4.88'
$ws.Range("E37").Value = '  +2.96%  '
$ws.Range("D38").Value = '5.65'
$ws.Range("E38").Value = '  +5.94%  '
$ws.Range("D39").Value = '155.01'
$ws.Range("E39").Value = '  +2.18%  '
$ws.Range("E40").Value = '  +4.66%  '
$ws.Range("D41").Value = '0.373'
$ws.Range("E41").Value = '  +1.37%  '
$ws.Range("D42").Value = '18.95'
$ws.Range("E42").Value = '  +1.89%  '
$ws.Range("D43").Value = '1.82'
$ws.Range("E43").Value = '  +5.73%  '
$ws.Range("D44").Value = '161.95'
$ws.Range("E44").Value = '  +5.72%  '
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("D46").Value = '0.0₆0302'
$ws.Range("E46").Value = '  -0.37%  '
$ws.Range("D47").Value = '15.72'
$ws.Range("E47").Value = '  +2.77%  '
$ws.Range("E48").Value = '  +3.00%  '
$ws.Range("D49").Value = '21.62'
$ws.Range("E49").Value = '  +6.22%  '
$ws.Range("D50").Value = '0.628'
$ws.Range("E50").Value = '  +3.81%  '
$ws.Range("D51").Value = '0.0520'
$ws.Range("E51").Value = '  +2.77%  '

# Restore the original (default/no explicit number format) styling so the
# cells' style index matches the source workbook (no s="..." attribute).
$ws.Range("D2:E51").ClearFormats()

